# Weekly update: insert one new price record at row 54, shifting the
# existing rows (and the old last row) down by one. Final used range
# grows from A1:R98 to A1:R99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54, pushing rows 54:98 down to 55:99.
# -4121 = xlShiftDown
$ws.Rows.Item(54).Insert(-4121)

# Populate the newly inserted row 54 with the new record.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = "2023-12-11"
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112028
$ws.Range("G54").Value = "Sandia"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 650
$ws.Range("K54").Value = 500
$ws.Range("L54").Value = 530
$ws.Range("M54").Value = 518
$ws.Range("N54").Value = "$/kilo (volumen en unidades)"
$ws.Range("O54").Value = "Perú"
$ws.Range("P54").Value = 518
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"
